# Update odds/value cells on rows 2-5 to the new figures, then remove the
# last match (row 6 - "Los Angeles FC" vs "Seattle Sounders") entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Macarthur FC - Auckland FC) ---
$ws.Range("H2").Value = 3.75
$ws.Range("Q2").Value = 1.67
$ws.Range("R2").Value = 2.2
$ws.Range("AC2").Value = 15
$ws.Range("AF2").Value = 41
$ws.Range("AJ2").Value = 9.5

# --- Row 3 (Daegu - Incheon) ---
$ws.Range("G3").Value = 2.6
$ws.Range("I3").Value = 2.6
$ws.Range("L3").Value = 3.4
$ws.Range("M3").Value = 1.07
$ws.Range("N3").Value = 9
$ws.Range("Z3").Value = 26
$ws.Range("AA3").Value = 21
$ws.Range("AH3").Value = 8.5
$ws.Range("AI3").Value = 13
$ws.Range("AN3").Value = 4.5
$ws.Range("BB3").Value = 81

# --- Row 4 (Daejeon - Jeju Utd) ---
$ws.Range("G4").Value = 2.2
$ws.Range("H4").Value = 3.2
$ws.Range("J4").Value = 3
$ws.Range("M4").Value = 1.06
$ws.Range("N4").Value = 10
$ws.Range("Y4").Value = 9.5
$ws.Range("AA4").Value = 19
$ws.Range("AB4").Value = 29
$ws.Range("AC4").Value = 9.5
$ws.Range("AD4").Value = 6
$ws.Range("AE4").Value = 13
$ws.Range("AI4").Value = 15
$ws.Range("AO4").Value = 13
$ws.Range("AP4").Value = 23
$ws.Range("AR4").Value = 67
$ws.Range("AY4").Value = 17

# --- Row 5 (Gwangju FC - Jeonbuk) ---
$ws.Range("G5").Value = 2.6
$ws.Range("I5").Value = 2.6
$ws.Range("J5").Value = 3.4
$ws.Range("M5").Value = 1.07
$ws.Range("N5").Value = 9
$ws.Range("O5").Value = 1.33
$ws.Range("P5").Value = 3.25
$ws.Range("X5").Value = 13
$ws.Range("AB5").Value = 34
$ws.Range("AC5").Value = 9
$ws.Range("AP5").Value = 26
$ws.Range("AS5").Value = 201
$ws.Range("AX5").Value = 4.5

# --- Remove row 6 (Los Angeles FC - Seattle Sounders) entirely ---
$ws.Rows(6).Delete()
